# Edit script for Lab_06_Bonus_Extra_help.docx
# Applies:
#  1. Text fix inside the "PPrint Module" explanation paragraph
#     ("check the the example for the weather.py)" -> "check the weather_example.py)")
#     reproduced with the exact target run layout.
#  2. A new "Tips" block (8 paragraphs) appended at the end of the document,
#     after the existing trailing empty paragraph.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Fix the "PPrint Module" paragraph text / run layout.
# ---------------------------------------------------------------------------

$pprintParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Short for Pretty Print*") {
        $pprintParaIndex = $i
        break
    }
}

if ($pprintParaIndex -eq -1) {
    throw "Could not locate the 'Short for Pretty Print' paragraph"
}

$pprintXml = '<w:p ' + $wNs + '><w:pPr><w:ind w:left="720"/></w:pPr>' +
    '<w:r><w:t>Short for Pretty Print, somewhat similar to print statements, but it prints out a nicer formatted form of data structures, for example printing nested dictionaries will be printed as a Left aligned tree, instead of one big block of code.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">  (</w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>check</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>' +
    '<w:r><w:t xml:space="preserve"> the </w:t></w:r>' +
    '<w:r><w:t>weather</w:t></w:r>' +
    '<w:r><w:t>_example</w:t></w:r>' +
    '<w:r><w:t>.py)</w:t></w:r>' +
    '</w:p>'

$d.Paragraphs.Item($pprintParaIndex).Range.InsertXML($pprintXml)

# ---------------------------------------------------------------------------
# 2. Append the new "Tips" block at the end of the document (after the
#    existing trailing blank paragraph, before the end of the body).
# ---------------------------------------------------------------------------

# Create 8 placeholder paragraphs after the current last (blank) paragraph.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
for ($n = 0; $n -lt 8; $n++) {
    $lastRange.InsertParagraphAfter()
}

# Index (1-based) of the first placeholder paragraph -- i.e. right after the
# paragraph that was previously last in the document.
$base = $lastPara.Index + 1

$newParaXml = @(
    # Para 1: "Tips: to use help() with custom modules, ..."
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="360"/></w:pPr>' +
        '<w:r><w:rPr><w:color w:val="FF0000"/><w:u w:val="single"/></w:rPr><w:t>Tips:</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
        '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t xml:space="preserve">to use </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>help(</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:rPr><w:i/><w:iCs/></w:rPr><w:t>) with custom modules, just run them once, and then  help(&lt;module name&gt;).</w:t></w:r>' +
        '</w:p>'),

    # Para 2: "A custom module made to easily get weather information from google weather."
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="360"/></w:pPr>' +
        '<w:r><w:t xml:space="preserve">A custom module made to easily get weather information from </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>google</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> weather.</w:t></w:r>' +
        '</w:p>'),

    # Para 3: "usage : " (underlined)
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="360"/><w:rPr><w:u w:val="single"/></w:rPr></w:pPr>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>usage :</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
        '</w:p>'),

    # Para 4: ">>> import weather" (bold)
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="1080"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">&gt;&gt;&gt; </w:t></w:r>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>import weather</w:t></w:r>' +
        '</w:p>'),

    # Para 5: ">>> data = weather.get_weather(Location)" (bold)
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="1080"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">&gt;&gt;&gt; </w:t></w:r>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>data</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>weather.get_weather</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>(Location)</w:t></w:r>' +
        '</w:p>'),

    # Para 6: empty (bold paragraph mark only)
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="1080"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'),

    # Para 7: "Location => A string with the location desired, e.g. "McMaster University""
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="1080"/></w:pPr>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Location</w:t></w:r>' +
        '<w:r><w:t xml:space="preserve"> =&gt; A string with the location desired, e.g. "McMaster University"</w:t></w:r>' +
        '</w:p>'),

    # Para 8: "data => output,is a dictionary"
    ('<w:p ' + $wNs + '><w:pPr><w:ind w:left="1080"/></w:pPr>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>data</w:t></w:r>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> =&gt; </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>output,is</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> a dictionary</w:t></w:r>' +
        '</w:p>')
)

for ($i = 0; $i -lt $newParaXml.Length; $i++) {
    $target = $d.Paragraphs.Item($base + $i)
    $target.Range.InsertXML($newParaXml[$i])
}

Write-Host "Done. Final paragraph count: $($d.Paragraphs.Count)"
